# chore: Adicionei uma coluna no Excel
# Adds an "updated_at" column (J) to Planilha1, between "created_at" (I) and
# "dispatch_number" (now K). Also fixes a couple of data values and the
# filtered range now that the sheet is one column wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- 1. Insert the new column -------------------------------------------
# Shifts former J:O -> K:P, copying formatting from the left neighbour
# (matches what the target file shows: J1 inherits I1's style, J2 inherits
# I2's style).
$ws.Columns.Item(10).Insert()

# New header + value for the inserted "updated_at" column.
$ws.Range("J1").Value = "updated_at"
$ws.Range("J2").Value = "2025-04-04 02:16:00.000 -0300"

# --- 2. Small data fixes on the untouched columns ------------------------
$ws.Range("C2").Value = 5
$ws.Range("H2").Value = "HPT-626041-1"

# --- 3. Row heights to match the refreshed layout -------------------------
$ws.Rows.Item(1).RowHeight = 36.6
$ws.Rows.Item(2).RowHeight = 36.6
$ws.Rows.Item(5).RowHeight = 52.8

# --- 4. Re-apply the AutoFilter over the new, wider range -----------------
$ws.AutoFilterMode = $false
$ws.Range("A1:P168").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$wb.Names.Item("Planilha1!_FilterDatabase").RefersTo = "=Planilha1!`$A`$1:`$P`$2"

# --- 5. Selection, matching the saved workbook's cursor position ----------
$ws.Range("J2").Select()
